# Auto-generated Excel COM-interop script
# Applies the scrape-refresh update described in the commit diff
# (new scrape batch at 14:17:13, 183/30/30 total rows across the 3 sheets)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Sheet: LP1912 ---
$ws1.Cells.Item(2, 1).Value = "Última actualización: 14:17:13"
$ws1.Cells.Item(3, 1).Value = "Total filas: 183"
$ws1.Cells.Item(15, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(16, 3).Value = "225_GOMEZ"
$ws1.Cells.Item(23, 1).Value = "06:46:40"
$ws1.Cells.Item(23, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(23, 4).Value = 35
$ws1.Cells.Item(24, 1).Value = "06:15:23"
$ws1.Cells.Item(24, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(24, 4).Value = 66
$ws1.Cells.Item(43, 1).Value = "06:58:58"
$ws1.Cells.Item(43, 3).Value = "15_ABASTO"
$ws1.Cells.Item(43, 4).Value = 90
$ws1.Cells.Item(44, 1).Value = "08:14:55"
$ws1.Cells.Item(44, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(44, 4).Value = 14
$ws1.Cells.Item(45, 1).Value = "08:14:55"
$ws1.Cells.Item(45, 3).Value = "15_ABASTO"
$ws1.Cells.Item(45, 4).Value = 15
$ws1.Cells.Item(46, 1).Value = "07:51:40"
$ws1.Cells.Item(46, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(46, 4).Value = 38
$ws1.Cells.Item(64, 1).Value = "08:49:06"
$ws1.Cells.Item(64, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(64, 4).Value = 28
$ws1.Cells.Item(65, 1).Value = "08:57:42"
$ws1.Cells.Item(65, 3).Value = "14_ABASTO"
$ws1.Cells.Item(65, 4).Value = 20
$ws1.Cells.Item(66, 1).Value = "08:49:06"
$ws1.Cells.Item(66, 3).Value = "14_ABASTO"
$ws1.Cells.Item(66, 4).Value = 29
$ws1.Cells.Item(67, 1).Value = "08:57:42"
$ws1.Cells.Item(67, 3).Value = "15X38_ABASTO"
$ws1.Cells.Item(67, 4).Value = 21
$ws1.Cells.Item(117, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(118, 3).Value = "17_ROMERO"
$ws1.Cells.Item(122, 1).Value = "11:56:32"
$ws1.Cells.Item(122, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(122, 4).Value = 14
$ws1.Cells.Item(123, 1).Value = "11:38:09"
$ws1.Cells.Item(123, 3).Value = "15_ABASTO"
$ws1.Cells.Item(123, 4).Value = 32
$ws1.Cells.Item(137, 3).Value = "15X38_ABASTO"
$ws1.Cells.Item(138, 3).Value = "14_ABASTO"
$ws1.Cells.Item(148, 1).Value = "12:58:23"
$ws1.Cells.Item(148, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(148, 4).Value = 27
$ws1.Cells.Item(149, 1).Value = "12:43:13"
$ws1.Cells.Item(149, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(149, 4).Value = 42
$ws1.Cells.Item(151, 1).Value = "12:43:13"
$ws1.Cells.Item(151, 3).Value = "14_ABASTO"
$ws1.Cells.Item(151, 4).Value = 50
$ws1.Cells.Item(152, 1).Value = "13:28:27"
$ws1.Cells.Item(152, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(152, 4).Value = 5
$ws1.Cells.Item(162, 1).Value = "14:17:13"
$ws1.Cells.Item(162, 4).Value = 0
$ws1.Cells.Item(163, 1).Value = "14:17:13"
$ws1.Cells.Item(163, 4).Value = 0
$ws1.Cells.Item(165, 1).Value = "14:17:13"
$ws1.Cells.Item(165, 4).Value = 10
$ws1.Cells.Item(166, 1).Value = "14:17:13"
$ws1.Cells.Item(166, 4).Value = 15
$ws1.Cells.Item(167, 1).Value = "14:17:13"
$ws1.Cells.Item(167, 4).Value = 17
$ws1.Cells.Item(168, 1).Value = "14:17:13"
$ws1.Cells.Item(168, 4).Value = 22
$ws1.Cells.Item(170, 1).Value = "14:17:13"
$ws1.Cells.Item(170, 2).Value = "14:48"
$ws1.Cells.Item(170, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(170, 4).Value = 31
$ws1.Cells.Item(171, 1).Value = "14:17:13"
$ws1.Cells.Item(171, 4).Value = 34
$ws1.Cells.Item(172, 1).Value = "14:17:13"
$ws1.Cells.Item(172, 2).Value = "14:51"
$ws1.Cells.Item(172, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(172, 4).Value = 34
$ws1.Cells.Item(173, 1).Value = "14:17:13"
$ws1.Cells.Item(173, 2).Value = "14:54"
$ws1.Cells.Item(173, 3).Value = "215_EL PELIGRO"
$ws1.Cells.Item(173, 4).Value = 37
$ws1.Cells.Item(174, 1).Value = "14:17:13"
$ws1.Cells.Item(174, 2).Value = "15:02"
$ws1.Cells.Item(174, 3).Value = "10_OLMOS"
$ws1.Cells.Item(174, 4).Value = 45
$ws1.Cells.Item(175, 1).Value = "13:28:27"
$ws1.Cells.Item(175, 2).Value = "15:11"
$ws1.Cells.Item(175, 3).Value = "14_ABASTO"
$ws1.Cells.Item(175, 4).Value = 103
$ws1.Cells.Item(176, 1).Value = "14:17:13"
$ws1.Cells.Item(176, 2).Value = "15:13"
$ws1.Cells.Item(176, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(176, 4).Value = 56
$ws1.Cells.Item(177, 2).Value = "15:17"
$ws1.Cells.Item(177, 3).Value = "14_ABASTO"
$ws1.Cells.Item(177, 4).Value = 83
$ws1.Cells.Item(178, 1).Value = "14:17:13"
$ws1.Cells.Item(178, 2).Value = "15:18"
$ws1.Cells.Item(178, 3).Value = "14_ABASTO"
$ws1.Cells.Item(178, 4).Value = 61
$ws1.Cells.Item(179, 1).Value = "14:17:13"
$ws1.Cells.Item(179, 2).Value = "15:29"
$ws1.Cells.Item(179, 3).Value = "10_OLMOS"
$ws1.Cells.Item(179, 4).Value = 72
$ws1.Cells.Item(180, 1).Value = "14:17:13"
$ws1.Cells.Item(180, 2).Value = "15:34"
$ws1.Cells.Item(180, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(180, 4).Value = 77
$ws1.Cells.Item(181, 1).Value = "14:17:13"
$ws1.Cells.Item(181, 2).Value = "15:36"
$ws1.Cells.Item(181, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(181, 4).Value = 79
$ws1.Cells.Item(181, 5).Value = "LP1912"
$ws1.Cells.Item(182, 1).Value = "13:54:35"
$ws1.Cells.Item(182, 2).Value = "15:41"
$ws1.Cells.Item(182, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(182, 4).Value = 107
$ws1.Cells.Item(182, 5).Value = "LP1912"
$ws1.Cells.Item(183, 1).Value = "14:17:13"
$ws1.Cells.Item(183, 2).Value = "15:42"
$ws1.Cells.Item(183, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(183, 4).Value = 85
$ws1.Cells.Item(183, 5).Value = "LP1912"
$ws1.Cells.Item(184, 1).Value = "14:17:13"
$ws1.Cells.Item(184, 2).Value = "15:53"
$ws1.Cells.Item(184, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(184, 4).Value = 96
$ws1.Cells.Item(184, 5).Value = "LP1912"
$ws1.Cells.Item(185, 1).Value = "14:17:13"
$ws1.Cells.Item(185, 2).Value = "15:53"
$ws1.Cells.Item(185, 3).Value = "15X38_ABASTO"
$ws1.Cells.Item(185, 4).Value = 96
$ws1.Cells.Item(185, 5).Value = "LP1912"
$ws1.Cells.Item(186, 1).Value = "14:17:13"
$ws1.Cells.Item(186, 2).Value = "15:57"
$ws1.Cells.Item(186, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(186, 4).Value = 100
$ws1.Cells.Item(186, 5).Value = "LP1912"
$ws1.Cells.Item(187, 1).Value = "14:17:13"
$ws1.Cells.Item(187, 2).Value = "16:06"
$ws1.Cells.Item(187, 3).Value = "14_ABASTO"
$ws1.Cells.Item(187, 4).Value = 109
$ws1.Cells.Item(187, 5).Value = "LP1912"
$ws1.Cells.Item(188, 1).Value = "14:17:13"
$ws1.Cells.Item(188, 2).Value = "16:14"
$ws1.Cells.Item(188, 3).Value = "17_ROMERO"
$ws1.Cells.Item(188, 4).Value = 117
$ws1.Cells.Item(188, 5).Value = "LP1912"

# --- Sheet: LP1912-215 ---
$ws2.Cells.Item(2, 1).Value = "Última actualización: 14:17:13"
$ws2.Cells.Item(3, 1).Value = "Total filas: 30"
$ws2.Cells.Item(31, 1).Value = "14:17:13"
$ws2.Cells.Item(31, 4).Value = 17
$ws2.Cells.Item(33, 1).Value = "14:17:13"
$ws2.Cells.Item(33, 2).Value = "14:48"
$ws2.Cells.Item(33, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(33, 4).Value = 31
$ws2.Cells.Item(34, 1).Value = "14:17:13"
$ws2.Cells.Item(34, 2).Value = "14:54"
$ws2.Cells.Item(34, 3).Value = "215_EL PELIGRO"
$ws2.Cells.Item(34, 4).Value = 37
$ws2.Cells.Item(35, 1).Value = "14:17:13"
$ws2.Cells.Item(35, 2).Value = "15:34"
$ws2.Cells.Item(35, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(35, 4).Value = 77
$ws2.Cells.Item(35, 5).Value = "LP1912"

# --- Sheet: 6203-6173 ---
$ws3.Cells.Item(2, 1).Value = "Última actualización: 14:17:13"
$ws3.Cells.Item(3, 1).Value = "Total filas: 30"
$ws3.Cells.Item(19, 1).Value = "08:49:06"
$ws3.Cells.Item(19, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(19, 4).Value = 101
$ws3.Cells.Item(20, 1).Value = "09:42:42"
$ws3.Cells.Item(20, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(20, 4).Value = 48
$ws3.Cells.Item(33, 1).Value = "14:17:13"
$ws3.Cells.Item(33, 4).Value = 10
$ws3.Cells.Item(34, 1).Value = "14:17:13"
$ws3.Cells.Item(34, 4).Value = 65
$ws3.Cells.Item(35, 1).Value = "14:17:13"
$ws3.Cells.Item(35, 2).Value = "16:02"
$ws3.Cells.Item(35, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(35, 4).Value = 105
$ws3.Cells.Item(35, 5).Value = "L6203"

